# Generate Report for Handoff
#
# The handoff report previously reflected a "handed back" state; this
# regenerates it for a fresh handoff: the status text moves from
# "Handed back: in sync with en-US" to "Ready for handoff", the
# handoff/generate timestamps are refreshed, and the now-shorter status
# text lets the status/date columns narrow on every sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps (left as plain text, matching source formatting)
$wsOverview.Range("G2").Value = "2016-12-15 05:03:53"
$wsDeDe.Range("H2").Value     = "2016-12-15 05:03:53"
$wsZhCn.Range("H2").Value     = "2016-12-15 05:03:39"

# --- The "Ready for handoff" status is shorter than the old text, so the
#     status/date columns can shrink to fit the new content.
$wsOverview.Columns("E:E").ColumnWidth = 16.33
$wsOverview.Columns("F:F").ColumnWidth = 16.33
$wsZhCn.Columns("C:C").ColumnWidth     = 16.33
$wsDeDe.Columns("C:C").ColumnWidth     = 16.33
